$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-504) holds the "Förändrad" (Changed) date as an Excel
# serial number. Every value of 45177 (2023-09-08) is bumped to 45178
# (2023-09-09).
$ws.Range("C2:C504").Value = 45178
